$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Insert-BodyXmlAt($range, [string]$bodyFragment) {
    $xml = $pkgHeader + $bodyFragment + $pkgFooter
    $null = $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 0. Drop the old `_GoBack` bookmark that currently sits at the very end of
#    the document - it gets recreated further down, right after the text
#    block we insert in step 1.
# ---------------------------------------------------------------------------
try {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
} catch { }

# ---------------------------------------------------------------------------
# 1. Insert nine new paragraphs right before the "Przechwycony z Monitorow"
#    paragraph (which currently immediately follows the "Skrocony opis
#    dzialania" heading paragraph).
# ---------------------------------------------------------------------------

$find1 = $d.Content
$null = $find1.Find.Execute("Przechwycony z ")
$insertPoint = $d.Range($find1.Start, $find1.Start)

$newParasFragment = '<w:p><w:r><w:t>Po uruchomieniu komponent próbuje pobrać z Monitor</w:t></w:r><w:r><w:t xml:space="preserve">ów dane zasobów i </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dostepnych</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> na</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t xml:space="preserve">tych zasobach pomiarów. Informacje te pobierane są za pomocą klasy </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>WebRequest</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Jeśli </w:t></w:r></w:p>' +
    '<w:p><w:r><w:t xml:space="preserve">krok ten się powiedzie, w oddzielnym wątku uruchamiane są operacje odpowiedzialne za </w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>nasłuchiwanie zapytań HTTP.</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t xml:space="preserve">Rejestrowanie zapytań HTTP realizowane jest za pomocą klasy </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>HttpListener</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> dostępnej </w:t></w:r></w:p>' +
    '<w:p><w:r><w:t xml:space="preserve">standardowo w przestrzeni nazw System.Net. W momencie wykrycia zapytania klasa ta </w:t></w:r></w:p>' +
    '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>przechwyca</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> wymagane dane i przekazuje je do metody odpowiedzialnej za sprawdzenie </w:t></w:r></w:p>' +
    '<w:p><w:r><w:t xml:space="preserve">poprawności ścieżki i (jeśli jest to wymagane) wygenerowania odpowiedzi, a następnie </w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>wraca do trybu nasłuchiwania.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
    '<w:p></w:p>'

Insert-BodyXmlAt $insertPoint $newParasFragment

# InsertXML's trailing (10th, empty) paragraph mark in the fragment above
# merges into the destination paragraph, leaving an extra empty paragraph
# behind it; delete that leftover so "Przechwycony z ..." stays untouched.
$find1b = $d.Content
$null = $find1b.Find.Execute("Przechwycony z ")
$destPara = $d.Range($find1b.Start, $find1b.Start).Paragraphs.First
$leftover = $destPara.Previous()
if ($leftover.Range.Text -eq "`r") {
    $leftover.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2. Move `lastRenderedPageBreak`: strip it from the run beginning
#    "W wypadku przyjsciu komunikatu z " and add it to the start of the
#    "Katalog" run that begins "Katalog przetwarza JSON...".
# ---------------------------------------------------------------------------
$fr = $d.Content
$null = $fr.Find.Execute("W wypadku przyjściu komunikatu z ")
$oldBreakRange = $d.Range($fr.Start, $fr.End)
$oldBreakRange.Delete()
$reinsertPt = $d.Range($fr.Start, $fr.Start)
Insert-BodyXmlAt $reinsertPt '<w:p><w:r><w:t xml:space="preserve">W wypadku przyjściu komunikatu z </w:t></w:r></w:p>'

$fk = $d.Content
$null = $fk.Find.Execute("Katalog przetwarza")
$katalogRunRange = $d.Range($fk.Start, $fk.Start + 7)
$katalogRunRange.Delete()
$katalogInsPt = $d.Range($fk.Start, $fk.Start)
Insert-BodyXmlAt $katalogInsPt '<w:p><w:r w:rsidRPr="00FD3956"><w:rPr><w:rStyle w:val="Wyrnienieintensywne"/></w:rPr><w:lastRenderedPageBreak/><w:t>Katalog</w:t></w:r></w:p>'

Write-Output "done"
